$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Turn the old two BodyText paragraphs describing the clinic example into:
#      - a new "Comparison to an overall mean" Heading3
#      - a FirstParagraph with merged/edited text
#      - a BodyText with the remaining original text
#    and wrap the three of them in a "comparison-to-an-overall-mean" bookmark.
# ---------------------------------------------------------------------------

# Insert a brand-new heading paragraph right before the old first BodyText
# paragraph of this section.
$pClinic = $d.Paragraphs(5)
$pClinic.Range.InsertParagraphBefore()
$newHeading = $d.Paragraphs(5)
$newHeading.Range.Text = "Comparison to an overall mean"
$newHeading.Style = "Heading3"

# The paragraph that used to start "You are in a setting..." becomes the new
# FirstParagraph; restyle it and rewrite its tail so it flows into the
# (trimmed) start of the next paragraph.
$pIntro = $d.Paragraphs(6)
$pIntro.Style = "FirstParagraph"

$d.Content.Find.Execute(
    "The patients themselves may differ markedly from one to another, but they are assigned in a more or less random way to each clinic. So, although one patient may differ markedly from another, the average result should be the same for across all the clinics.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The time spent with each patient, the total cost of the medications prescribed, the length of the notes written in the medical record–these can differ from patient to patient. But because of how the patients are assigned to a particular clinic (more or less randomly), these should balance out on average.",
    2) | Out-Null

# Remove the now-duplicated lead-in sentence from the following paragraph
# (leaves "If there is one clinic..." onward, including the embedded quote
# runs, untouched).
$d.Content.Find.Execute(
    "The time spent with each patient, the total cost of the medications prescribed, the length of the notes written in the medical record–these can differ from patient to patient. But because of how the patients are assigned to a particular clinic, these should balance out on average. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2) | Out-Null

# Wrap the new heading + the two body paragraphs with the section bookmark.
$bmRange = $d.Range($d.Paragraphs(5).Range.Start, $d.Paragraphs(7).Range.End)
$d.Bookmarks.Add("comparison-to-an-overall-mean", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Split the "Suppose you are in a setting..." paragraph: the figure moves
#    into its own new BodyText paragraph, and the stray joining space is
#    removed.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Then you would write your hypothesis as") | Out-Null
$endOfSentence = $findRng.End

$hyp = $d.Content
$hyp.Find.Execute("Then you would write your hypothesis as") | Out-Null
$hypPara = $hyp.Paragraphs(1)
$splitPoint = $hypPara.Range.Duplicate
$splitPoint.SetRange($endOfSentence, $endOfSentence)
$splitPoint.InsertParagraphAfter()

# Locate the freshly split-off paragraph (now holds just the leading space
# run and the drawing) and clean it up.
$figPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs($i).Range.Text
    if ($ptxt -like "*Then you would write your hypothesis as*") {
        $figPara = $d.Paragraphs($i + 1)
        break
    }
}
$figPara.Style = "BodyText"
$leadSpace = $figPara.Range.Duplicate
$leadSpace.SetRange($figPara.Range.Start, $figPara.Range.Start + 1)
$leadSpace.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Plain text edits.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The analysis of means approach that compares each group mean to the overall mean is easy to implement and it lends itself to a simple graphical display.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The analysis of means approach compares each group mean to the overall mean. It is easy to implement and it lends itself to a simple graphical display. You need a table of critical values, which depend on alpha (the overall Type I error rate), g (the number of groups), and n (the number of observations within each group). Some tables use the degrees of freedom for error in place of n.",
    2) | Out-Null

$d.Content.Find.Execute(
    "You are using statistics to help in understanding why deviations from the norm occur.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You are using statistics to help in understanding if deviations from the norm occur and then study those deviating groups to understand why they deviate.",
    2) | Out-Null

$d.Content.Find.Execute(
    "The research question is whether all strains have a comparable amount of CryA1c in their leaves.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The research question is whether any strain differs from the overall mean level of CryA1c.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) New paragraph about one-sided tests, right after the "Notice that all
#    five means..." paragraph (and therefore still before the existing
#    "a-simple-example" bookmark end).
# ---------------------------------------------------------------------------
$noticeRng = $d.Content
$noticeRng.Find.Execute("None of the five strains shows a statisticially significant difference from the overall mean.") | Out-Null
$noticePara = $noticeRng.Paragraphs(1)
$noticePara.Range.InsertParagraphAfter()

$oneSidedPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Notice that all five means lie inside the limits*") {
        $oneSidedPara = $d.Paragraphs($i + 1)
        break
    }
}
$oneSidedPara.Range.Text = "You may prefer a one-sided test in this setting, such as testing whether any strain is deficient in the CryA1c levels. You can make a very easy modification to get one-sided tests."
$oneSidedPara.Style = "BodyText"

# ---------------------------------------------------------------------------
# 5) New paragraph about specifying the hypothesis before looking at data,
#    inserted right before the closing "The Analysis of Means, just like..."
#    paragraph (still before the existing "caveats" bookmark end).
# ---------------------------------------------------------------------------
$closingRng = $d.Content
$closingRng.Find.Execute("The Analysis of Means, just like the Dunnett") | Out-Null
$closingPara = $closingRng.Paragraphs(1)
$closingPara.Range.InsertParagraphBefore()

$peekingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*The Analysis of Means, just like*") {
        $peekingPara = $d.Paragraphs($i - 1)
        break
    }
}
$peekingPara.Range.Text = "You also need to specify the Analysis of Means hypothesis prior to looking at your data. Peeking at the data and then choosing your hypothesis is cheating."
$peekingPara.Style = "BodyText"

Write-Host "Edit complete."
